$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.958.43"
$ws.Range("E2").Value = "  +1.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.622.73"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.06"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.05"
$ws.Range("E6").Value = "  +2.54%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  +0.30%  "

$ws.Range("E9").Value = "  +1.24%  "

$ws.Range("E10").Value = "  +5.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.61"
$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.151"
$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.56"
$ws.Range("E13").Value = "  +1.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.090.28"
$ws.Range("E14").Value = "  -0.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.802.70"
$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000149"
$ws.Range("E16").Value = "  +2.74%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.621.60"
$ws.Range("E17").Value = "  -0.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.18"
$ws.Range("E18").Value = "  +6.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.63"
$ws.Range("E19").Value = "  +3.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "349.83"
$ws.Range("E20").Value = "  +2.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.89"
$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.70"
$ws.Range("E23").Value = "  +2.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.27"
$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.73"
$ws.Range("E25").Value = "  +14.05%  "

$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.69"
$ws.Range("E26").Value = "  +1.80%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.21"
$ws.Range("E27").Value = "  +5.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.12"
$ws.Range("E28").Value = "  +3.82%  "

$ws.Range("E29").Value = "  +0.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "545.13"
$ws.Range("E30").Value = "  -1.79%  "

$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("E32").Value = "  +1.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0848"
$ws.Range("E33").Value = "  +5.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.76"
$ws.Range("E34").Value = "  +0.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.24"
$ws.Range("E35").Value = "  +0.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.01"
$ws.Range("E36").Value = "  +1.17%  "

$ws.Range("E37").Value = "  +0.42%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("E39").Value = "  +4.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.42"
$ws.Range("E40").Value = "  +2.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "168.29"
$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.89"
$ws.Range("E43").Value = "  +0.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.92"
$ws.Range("E44").Value = "  +4.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0588"
$ws.Range("E45").Value = "  +2.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.44"
$ws.Range("E46").Value = "  -4.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.627"
$ws.Range("E47").Value = "  +0.52%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0245"
$ws.Range("E48").Value = "  +0.81%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.99"
$ws.Range("E49").Value = "  +12.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0965"
$ws.Range("E50").Value = "  +0.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.15"
$ws.Range("E51").Value = "  +2.35%  "

